$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 3000  # ALC!H62 was 2502
$ws.Cells.Item(62, 9).Value = 3000  # ALC!I62 was 2502
$ws.Cells.Item(62, 11).Value = 3000  # ALC!K62 was 2502
$ws.Cells.Item(62, 13).Value = -2376  # ALC!M62 was -1878
$ws.Cells.Item(65, 8).Value = 3000  # ALC!H65 was 2502
$ws.Cells.Item(65, 9).Value = 3000  # ALC!I65 was 2502
$ws.Cells.Item(65, 11).Value = 15000  # ALC!K65 was 12510
$ws.Cells.Item(65, 13).Value = -11880  # ALC!M65 was -9390
$ws.Cells.Item(93, 8).Value = 29999  # ALC!H93 was 29999.5
$ws.Cells.Item(93, 10).Value = 29999  # ALC!J93 was 29999.5
$ws.Cells.Item(93, 12).Value = 29999  # ALC!L93 was 29999.5
$ws.Cells.Item(93, 14).Value = -34991  # ALC!N93 was -34991.5
$ws.Cells.Item(98, 8).Value = 1261.9642  # ALC!H98 was 1290.7407
$ws.Cells.Item(98, 9).Value = 1282.1538  # ALC!I98 was 1314.04
$ws.Cells.Item(98, 11).Value = 1282.1538  # ALC!K98 was 1314.04
$ws.Cells.Item(98, 13).Value = 215.8462  # ALC!M98 was 183.96
$ws.Cells.Item(106, 8).Value = 76957464  # ALC!H106 was 90942000
$ws.Cells.Item(106, 9).Value = 83364330  # ALC!I106 was 90942000
$ws.Cells.Item(106, 10).Value = 75000  # ALC!J106 was 0
$ws.Cells.Item(106, 11).Value = 83364330  # ALC!K106 was 90942000
$ws.Cells.Item(106, 12).Value = 75000  # ALC!L106 was 0
$ws.Cells.Item(106, 13).Value = -83363699  # ALC!M106 was -90941369
$ws.Cells.Item(106, 14).Value = -76262  # ALC!N106 was None
$ws.Cells.Item(122, 8).Value = 1261.9642  # ALC!H122 was 1290.7407
$ws.Cells.Item(122, 9).Value = 1282.1538  # ALC!I122 was 1314.04
$ws.Cells.Item(122, 11).Value = 3846.4614  # ALC!K122 was 3942.12
$ws.Cells.Item(122, 13).Value = -1396.4614  # ALC!M122 was -1492.12
$ws.Cells.Item(137, 8).Value = 2166.6667  # ALC!H137 was 1900
$ws.Cells.Item(137, 9).Value = 2000  # ALC!I137 was 1550
$ws.Cells.Item(137, 11).Value = 6000  # ALC!K137 was 4650
$ws.Cells.Item(137, 13).Value = -3450  # ALC!M137 was -2100
$ws.Cells.Item(138, 8).Value = 2429.9473  # ALC!H138 was 2506.611
$ws.Cells.Item(138, 9).Value = 1938.4445  # ALC!I138 was 1978.3846
$ws.Cells.Item(138, 10).Value = 3636.3635  # ALC!J138 was 3880
$ws.Cells.Item(138, 11).Value = 5815.333500000001  # ALC!K138 was 5935.1538
$ws.Cells.Item(138, 12).Value = 10909.0905  # ALC!L138 was 11640
$ws.Cells.Item(138, 13).Value = -675.3335000000006  # ALC!M138 was -795.1538
$ws.Cells.Item(138, 14).Value = -21189.0905  # ALC!N138 was -21920
$ws.Cells.Item(141, 8).Value = 3167.7812  # ALC!H141 was 3351.5334
$ws.Cells.Item(141, 9).Value = 1481.7142  # ALC!I141 was 1564.0385
$ws.Cells.Item(141, 11).Value = 4445.142599999999  # ALC!K141 was 4692.1155
$ws.Cells.Item(141, 13).Value = 734.8574000000008  # ALC!M141 was 487.8845000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 22223404  # ARM!H2 was 23810774
$ws.Cells.Item(2, 9).Value = 25641636  # ARM!I2 was 27778420
$ws.Cells.Item(2, 11).Value = 25641636  # ARM!K2 was 27778420
$ws.Cells.Item(2, 13).Value = -25641523  # ARM!M2 was -27778307
$ws.Cells.Item(23, 8).Value = 18500  # ARM!H23 was 0
$ws.Cells.Item(23, 10).Value = 18500  # ARM!J23 was 0
$ws.Cells.Item(23, 12).Value = 18500  # ARM!L23 was 0
$ws.Cells.Item(23, 14).Value = -19018  # ARM!N23 was None
$ws.Cells.Item(74, 8).Value = 526.1852  # ARM!H74 was 528.03705
$ws.Cells.Item(74, 9).Value = 526.1852  # ARM!I74 was 528.03705
$ws.Cells.Item(74, 11).Value = 526.1852  # ARM!K74 was 528.03705
$ws.Cells.Item(74, 13).Value = 347.8148  # ARM!M74 was 345.96295
$ws.Cells.Item(77, 8).Value = 526.1852  # ARM!H77 was 528.03705
$ws.Cells.Item(77, 9).Value = 526.1852  # ARM!I77 was 528.03705
$ws.Cells.Item(77, 11).Value = 2630.926  # ARM!K77 was 2640.18525
$ws.Cells.Item(77, 13).Value = 1737.074  # ARM!M77 was 1727.81475
$ws.Cells.Item(94, 8).Value = 0  # ARM!H94 was 65000
$ws.Cells.Item(94, 10).Value = 0  # ARM!J94 was 65000
$ws.Cells.Item(94, 12).Value = 0  # ARM!L94 was 65000
$ws.Cells.Item(94, 14).ClearContents()  # ARM!N94 was -66802
$ws.Cells.Item(110, 8).Value = 5052269.5  # ARM!H110 was 5292796.5
$ws.Cells.Item(110, 10).Value = 1249.75  # ARM!J110 was 1266.6666
$ws.Cells.Item(110, 12).Value = 1249.75  # ARM!L110 was 1266.6666
$ws.Cells.Item(110, 14).Value = -5339.75  # ARM!N110 was -5356.6666
$ws.Cells.Item(116, 8).Value = 22223404  # ARM!H116 was 23810774
$ws.Cells.Item(116, 9).Value = 25641636  # ARM!I116 was 27778420
$ws.Cells.Item(116, 11).Value = 25641636  # ARM!K116 was 27778420
$ws.Cells.Item(116, 13).Value = -25639342  # ARM!M116 was -27776126
$ws.Cells.Item(122, 8).Value = 1135190.6  # ARM!H122 was 1140373.4
$ws.Cells.Item(122, 9).Value = 1451316.6  # ARM!I122 was 1457980
$ws.Cells.Item(122, 11).Value = 4353949.800000001  # ARM!K122 was 4373940
$ws.Cells.Item(122, 13).Value = -4351499.800000001  # ARM!M122 was -4371490
$ws.Cells.Item(124, 8).Value = 75000  # ARM!H124 was 52714.5
$ws.Cells.Item(124, 10).Value = 75000  # ARM!J124 was 52714.5
$ws.Cells.Item(124, 12).Value = 75000  # ARM!L124 was 52714.5
$ws.Cells.Item(124, 14).Value = -84820  # ARM!N124 was -62534.5
$ws.Cells.Item(125, 8).Value = 0  # ARM!H125 was 75000
$ws.Cells.Item(125, 10).Value = 0  # ARM!J125 was 75000
$ws.Cells.Item(125, 12).Value = 0  # ARM!L125 was 75000
$ws.Cells.Item(125, 14).ClearContents()  # ARM!N125 was -84840

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 22223404  # BSM!H3 was 23810774
$ws.Cells.Item(3, 9).Value = 25641636  # BSM!I3 was 27778420
$ws.Cells.Item(3, 11).Value = 25641636  # BSM!K3 was 27778420
$ws.Cells.Item(3, 13).Value = -25641522  # BSM!M3 was -27778306
$ws.Cells.Item(20, 8).Value = 1184.2941  # BSM!H20 was 1141.1666
$ws.Cells.Item(20, 9).Value = 953.875  # BSM!I20 was 893.2222
$ws.Cells.Item(20, 11).Value = 953.875  # BSM!K20 was 893.2222
$ws.Cells.Item(20, 13).Value = -706.875  # BSM!M20 was -646.2222
$ws.Cells.Item(29, 8).Value = 1000  # BSM!H29 was 325
$ws.Cells.Item(29, 9).Value = 1000  # BSM!I29 was 325
$ws.Cells.Item(29, 11).Value = 1000  # BSM!K29 was 325
$ws.Cells.Item(29, 13).Value = -711  # BSM!M29 was -36
$ws.Cells.Item(107, 8).Value = 1559.5625  # BSM!H107 was 1597.2
$ws.Cells.Item(107, 9).Value = 1542.6923  # BSM!I107 was 1543.0769
$ws.Cells.Item(107, 10).Value = 1632.6666  # BSM!J107 was 1949
$ws.Cells.Item(107, 11).Value = 1542.6923  # BSM!K107 was 1543.0769
$ws.Cells.Item(107, 12).Value = 1632.6666  # BSM!L107 was 1949
$ws.Cells.Item(107, 13).Value = 377.3077000000001  # BSM!M107 was 376.9231
$ws.Cells.Item(107, 14).Value = -5472.6666  # BSM!N107 was -5789

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3664.0833  # CRP!H31 was 3784.261
$ws.Cells.Item(31, 9).Value = 2442.111  # CRP!I31 was 2532.8235
$ws.Cells.Item(31, 11).Value = 2442.111  # CRP!K31 was 2532.8235
$ws.Cells.Item(31, 13).Value = -2147.111  # CRP!M31 was -2237.8235
$ws.Cells.Item(34, 8).Value = 3664.0833  # CRP!H34 was 3784.261
$ws.Cells.Item(34, 9).Value = 2442.111  # CRP!I34 was 2532.8235
$ws.Cells.Item(34, 11).Value = 2442.111  # CRP!K34 was 2532.8235
$ws.Cells.Item(34, 13).Value = -2240.111  # CRP!M34 was -2330.8235

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 679  # CUL!H5 was 861.125
$ws.Cells.Item(5, 9).Value = 679.3333  # CUL!I5 was 922.5
$ws.Cells.Item(5, 10).Value = 678.6  # CUL!J5 was 799.75
$ws.Cells.Item(5, 11).Value = 2037.9999  # CUL!K5 was 2767.5
$ws.Cells.Item(5, 12).Value = 2035.8  # CUL!L5 was 2399.25
$ws.Cells.Item(5, 13).Value = -1925.9999  # CUL!M5 was -2655.5
$ws.Cells.Item(5, 14).Value = -2259.8  # CUL!N5 was -2623.25
$ws.Cells.Item(18, 8).Value = 704.8570999999999  # CUL!H18 was 752.5714
$ws.Cells.Item(18, 9).Value = 704.8570999999999  # CUL!I18 was 711.3333
$ws.Cells.Item(18, 10).Value = 0  # CUL!J18 was 1000
$ws.Cells.Item(18, 11).Value = 2114.5713  # CUL!K18 was 2133.9999
$ws.Cells.Item(18, 12).Value = 0  # CUL!L18 was 3000
$ws.Cells.Item(18, 13).Value = -1945.5713  # CUL!M18 was -1964.9999
$ws.Cells.Item(18, 14).ClearContents()  # CUL!N18 was -3338
$ws.Cells.Item(23, 8).Value = 333699.66  # CUL!H23 was 91024.45
$ws.Cells.Item(23, 9).Value = 0  # CUL!I23 was 15
$ws.Cells.Item(23, 10).Value = 333699.66  # CUL!J23 was 333716.34
$ws.Cells.Item(23, 11).Value = 0  # CUL!K23 was 45
$ws.Cells.Item(23, 12).Value = 1001098.98  # CUL!L23 was 1001149.02
$ws.Cells.Item(23, 13).ClearContents()  # CUL!M23 was 190
$ws.Cells.Item(23, 14).Value = -1001568.98  # CUL!N23 was -1001619.02
$ws.Cells.Item(81, 8).Value = 3400  # CUL!H81 was 3500
$ws.Cells.Item(81, 10).Value = 3400  # CUL!J81 was 3500
$ws.Cells.Item(81, 12).Value = 10200  # CUL!L81 was 10500
$ws.Cells.Item(81, 14).Value = -12446  # CUL!N81 was -12746
$ws.Cells.Item(84, 8).Value = 3400  # CUL!H84 was 3500
$ws.Cells.Item(84, 10).Value = 3400  # CUL!J84 was 3500
$ws.Cells.Item(84, 12).Value = 30600  # CUL!L84 was 31500
$ws.Cells.Item(84, 14).Value = -41832  # CUL!N84 was -42732
$ws.Cells.Item(115, 8).Value = 400  # CUL!H115 was 0
$ws.Cells.Item(115, 9).Value = 400  # CUL!I115 was 0
$ws.Cells.Item(115, 11).Value = 1200  # CUL!K115 was 0
$ws.Cells.Item(115, 13).Value = -25  # CUL!M115 was None
$ws.Cells.Item(117, 8).Value = 1741.5714  # CUL!H117 was 1090.125
$ws.Cells.Item(117, 9).Value = 0  # CUL!I117 was 699
$ws.Cells.Item(117, 10).Value = 1741.5714  # CUL!J117 was 1146
$ws.Cells.Item(117, 11).Value = 0  # CUL!K117 was 2097
$ws.Cells.Item(117, 12).Value = 5224.7142  # CUL!L117 was 3438
$ws.Cells.Item(117, 13).ClearContents()  # CUL!M117 was 1345
$ws.Cells.Item(117, 14).Value = -12108.7142  # CUL!N117 was -10322
$ws.Cells.Item(135, 8).Value = 679  # CUL!H135 was 861.125
$ws.Cells.Item(135, 9).Value = 679.3333  # CUL!I135 was 922.5
$ws.Cells.Item(135, 10).Value = 678.6  # CUL!J135 was 799.75
$ws.Cells.Item(135, 11).Value = 6113.9997  # CUL!K135 was 8302.5
$ws.Cells.Item(135, 12).Value = 6107.400000000001  # CUL!L135 was 7197.75
$ws.Cells.Item(135, 13).Value = -3578.9997  # CUL!M135 was -5767.5
$ws.Cells.Item(135, 14).Value = -11177.4  # CUL!N135 was -12267.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 34970.332  # GSM!H123 was 34979.89
$ws.Cells.Item(123, 9).Value = 34900  # GSM!I123 was 0
$ws.Cells.Item(123, 10).Value = 34984.4  # GSM!J123 was 34979.89
$ws.Cells.Item(123, 11).Value = 34900  # GSM!K123 was 0
$ws.Cells.Item(123, 12).Value = 34984.4  # GSM!L123 was 34979.89
$ws.Cells.Item(123, 13).Value = -32450  # GSM!M123 was None
$ws.Cells.Item(123, 14).Value = -39884.4  # GSM!N123 was -39879.89

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2844  # LTW!H40 was 2771.75
$ws.Cells.Item(40, 9).Value = 2862.3333  # LTW!I40 was 2771.75
$ws.Cells.Item(40, 10).Value = 2789  # LTW!J40 was 0
$ws.Cells.Item(40, 11).Value = 2862.3333  # LTW!K40 was 2771.75
$ws.Cells.Item(40, 12).Value = 2789  # LTW!L40 was 0
$ws.Cells.Item(40, 13).Value = -2726.3333  # LTW!M40 was -2635.75
$ws.Cells.Item(40, 14).Value = -3061  # LTW!N40 was None
$ws.Cells.Item(43, 8).Value = 339567.2  # LTW!H43 was 318400.5
$ws.Cells.Item(43, 9).Value = 4670  # LTW!I43 was 3727.5
$ws.Cells.Item(43, 11).Value = 4670  # LTW!K43 was 3727.5
$ws.Cells.Item(43, 13).Value = -4477  # LTW!M43 was -3534.5
$ws.Cells.Item(61, 8).Value = 37038704  # LTW!H61 was 27779778
$ws.Cells.Item(113, 8).Value = 37038704  # LTW!H113 was 27779778
$ws.Cells.Item(122, 8).Value = 0  # LTW!H122 was 3000
$ws.Cells.Item(122, 10).Value = 0  # LTW!J122 was 3000
$ws.Cells.Item(122, 12).Value = 0  # LTW!L122 was 9000
$ws.Cells.Item(122, 14).ClearContents()  # LTW!N122 was -13900
$ws.Cells.Item(132, 8).Value = 3000  # LTW!H132 was 950
$ws.Cells.Item(132, 9).Value = 3000  # LTW!I132 was 950
$ws.Cells.Item(132, 11).Value = 9000  # LTW!K132 was 2850
$ws.Cells.Item(132, 13).Value = -6470  # LTW!M132 was -320

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 2066.6667  # WVR!H9 was 4050
$ws.Cells.Item(9, 10).Value = 2350  # WVR!J9 was 4900
$ws.Cells.Item(9, 12).Value = 2350  # WVR!L9 was 4900
$ws.Cells.Item(9, 14).Value = -2630  # WVR!N9 was -5180
$ws.Cells.Item(122, 8).Value = 2629.8  # WVR!H122 was 2674.75
$ws.Cells.Item(122, 9).Value = 2629.8  # WVR!I122 was 2674.75
$ws.Cells.Item(122, 11).Value = 7889.400000000001  # WVR!K122 was 8024.25
$ws.Cells.Item(122, 13).Value = -5439.400000000001  # WVR!M122 was -5574.25
$ws.Cells.Item(132, 8).Value = 2683.2632  # WVR!H132 was 2846
$ws.Cells.Item(132, 9).Value = 2092.6875  # WVR!I132 was 2205.9285
$ws.Cells.Item(132, 11).Value = 6278.0625  # WVR!K132 was 6617.7855
$ws.Cells.Item(132, 13).Value = -3748.0625  # WVR!M132 was -4087.7855
